$d = $word.ActiveDocument

$replacements = @(
    @{old="846÷9="; new="495÷7="},
    @{old="115÷3="; new="368÷8="},
    @{old="319÷4="; new="775÷7="},
    @{old="705÷6="; new="807÷7="},
    @{old="993÷6="; new="781÷8="},
    @{old="183÷3="; new="780÷8="},
    @{old="989÷5="; new="673÷4="},
    @{old="832÷9="; new="512÷9="},
    @{old="837÷4="; new="675÷2="},
    @{old="288÷3="; new="759÷2="},
    @{old="644÷8="; new="211÷9="},
    @{old="655÷5="; new="502÷8="},
    @{old="148÷2="; new="821÷8="},
    @{old="142÷2="; new="172÷2="},
    @{old="800÷5="; new="784÷7="},
    @{old="819÷9="; new="251÷5="},
    @{old="944÷4="; new="811÷4="},
    @{old="556÷4="; new="953÷9="},
    @{old="963÷7="; new="623÷7="},
    @{old="262÷5="; new="350÷7="},
    @{old="881÷8="; new="812÷7="},
    @{old="109÷2="; new="245÷4="},
    @{old="711÷5="; new="703÷7="},
    @{old="576÷7="; new="622÷5="},
    @{old="846÷7="; new="480÷6="}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $pair.new, 2)
}
